# The presentation ships with two DrawingML theme parts:
#   ppt/theme/theme1.xml -> used by the Notes Master (originally "Office Theme")
#   ppt/theme/theme2.xml -> used by the Slide Master  (originally "Integral")
#
# The authored edit swaps the two themes' contents: the Slide Master now
# carries the plain "Office Theme" colour scheme, and the Notes Master now
# carries the "Integral" colour scheme. Font scheme / format scheme are
# identical between the two themes already, so only the 12 theme colour
# slots (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) need to move.

$p  = $ppt.ActivePresentation
$sm = $p.SlideMaster
$nm = $p.NotesMaster

$slideMasterColors = $sm.Theme.ThemeColorScheme
$notesMasterColors = $nm.Theme.ThemeColorScheme

# Target palette for the Slide Master's theme part (plain "Office Theme").
$slideMasterColors.Colors(1).RGB  = 0          # dk1      000000
$slideMasterColors.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$slideMasterColors.Colors(3).RGB  = 6968388    # dk2      44546A
$slideMasterColors.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$slideMasterColors.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$slideMasterColors.Colors(6).RGB  = 3243501    # accent2  ED7D31
$slideMasterColors.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$slideMasterColors.Colors(8).RGB  = 49407      # accent4  FFC000
$slideMasterColors.Colors(9).RGB  = 12874308   # accent5  4472C4
$slideMasterColors.Colors(10).RGB = 4697456    # accent6  70AD47
$slideMasterColors.Colors(11).RGB = 12673797   # hlink    0563C1
$slideMasterColors.Colors(12).RGB = 7491477    # folHlink 954F72

# Target palette for the Notes Master's theme part ("Integral").
$notesMasterColors.Colors(1).RGB  = 0          # dk1      000000
$notesMasterColors.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$notesMasterColors.Colors(3).RGB  = 5332805    # dk2      455F51
$notesMasterColors.Colors(4).RGB  = 13754083   # lt2      E3DED1
$notesMasterColors.Colors(5).RGB  = 3722137    # accent1  99CB38
$notesMasterColors.Colors(6).RGB  = 3646819    # accent2  63A537
$notesMasterColors.Colors(7).RGB  = 2412774    # accent3  E6D024
$notesMasterColors.Colors(8).RGB  = 38860      # accent4  CC9700
$notesMasterColors.Colors(9).RGB  = 13611854   # accent5  4EB3CF
$notesMasterColors.Colors(10).RGB = 10915127   # accent6  378DA6
$notesMasterColors.Colors(11).RGB = 2465643    # hlink    6B9F25
$notesMasterColors.Colors(12).RGB = 158642     # folHlink B26B02
